# The workbook is a small "quadratic-svm-score" style report with columns:
#   A: Row (sample ids)      B: 1-f__Clostridiaceae (scores)   C: prediction
#
# This commit refreshes the numeric score column (B) with newly computed
# values for the three data rows, while keeping the text labels / layout
# (and their text-formatted style) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the text formatting on the header row and the row-label column,
# matching how the previous export pass stamped these cells.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A4").NumberFormat = "@"

# Updated score values produced by the latest pipeline run.
$ws.Range("B2").Value = 1116.4816852553504
$ws.Range("B3").Value = 11634.999854082726
$ws.Range("B4").Value = 2263.3332204411518
